$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 407 (shifts existing rows 407:509 down to 408:510,
# dimension grows from A1:R509 to A1:R510).
$ws.Rows("407:407").Insert()

$ws.Cells.Item(407,1).Value = 9
$ws.Cells.Item(407,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(407,3).Value = "Metropolitana"
$ws.Cells.Item(407,4).Value = 44932
$ws.Cells.Item(407,5).Value = 13
$ws.Cells.Item(407,6).Value = 100112039
$ws.Cells.Item(407,7).Value = "Ciboulette"
$ws.Cells.Item(407,8).Value = "Sin especificar"
$ws.Cells.Item(407,9).Value = "Primera"
$ws.Cells.Item(407,10).Value = 340
$ws.Cells.Item(407,11).Value = 1000
$ws.Cells.Item(407,12).Value = 1000
$ws.Cells.Item(407,13).Value = 1000
$ws.Cells.Item(407,14).Value = "`$/docena de atados"
$ws.Cells.Item(407,15).Value = "Región Metropolitana"
$ws.Cells.Item(407,16).Value = 333
$ws.Cells.Item(407,17).Value = 3
$ws.Cells.Item(407,18).Value = "Hortaliza"
